$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the neighboring header cell (G1) onto the new
# header cell (H1) so the new column matches the existing header style
# (bold, centered, bordered), then set its text.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Add data values in the new "Save" column
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
